$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2024-12-17 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-18 Wednesday", 2)

# Update each "NN x NN=" multiplication problem in the table.
# Cell-targeted Range.Text assignment is used (instead of a global
# Find/Replace) because several of the old/new values collide with each
# other across different cells (e.g. a cell's new value equals another
# cell's old value), which would cause a document-wide replace to hit the
# wrong cell a second time.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "91×35="   # was 96×64=
$t.Cell(1,2).Range.Text  = "62×74="   # was 85×41=
$t.Cell(1,3).Range.Text  = "68×42="   # was 17×54=
$t.Cell(1,4).Range.Text  = "67×47="   # was 45×18=
$t.Cell(1,5).Range.Text  = "26×18="   # was 17×59=

$t.Cell(5,1).Range.Text  = "24×48="   # was 88×33=
$t.Cell(5,2).Range.Text  = "74×63="   # was 30×42=
$t.Cell(5,3).Range.Text  = "55×59="   # was 77×29=
$t.Cell(5,4).Range.Text  = "97×90="   # was 71×62=
$t.Cell(5,5).Range.Text  = "79×80="   # was 30×65=

$t.Cell(10,1).Range.Text = "81×67="   # was 96×90=
$t.Cell(10,2).Range.Text = "74×43="   # was 23×39=
$t.Cell(10,3).Range.Text = "75×57="   # was 58×40=
$t.Cell(10,4).Range.Text = "82×30="   # was 14×31=
$t.Cell(10,5).Range.Text = "27×55="   # was 87×96=

$t.Cell(15,1).Range.Text = "99×54="   # was 83×19=
$t.Cell(15,2).Range.Text = "98×35="   # was 18×75=
$t.Cell(15,3).Range.Text = "90×79="   # was 12×82=
$t.Cell(15,4).Range.Text = "30×27="   # was 79×75=
$t.Cell(15,5).Range.Text = "69×38="   # was 44×21=

$t.Cell(20,1).Range.Text = "15×39="   # was 26×68=
$t.Cell(20,2).Range.Text = "92×92="   # was 88×84=
$t.Cell(20,3).Range.Text = "40×13="   # was 66×22=
$t.Cell(20,4).Range.Text = "33×32="   # was 91×35=
$t.Cell(20,5).Range.Text = "29×51="   # was 74×45=
